$wb = $excel.ActiveWorkbook
$wsIn = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

# 1. Fix product name: add hyphen after "775"
$wsIn.Range("B1").Value = "775-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-Late Repayment"

# 2. Description cell now mirrors the short name "775C"
$wsIn.Range("B3").Value = "775C"
$wsIn.Range("B2").Copy()
$wsIn.Range("B3").PasteSpecial(-4122)

# 3. Currency label lower-cased, currency value trimmed of trailing space
$wsIn.Range("A6").Value = "currency"
$wsIn.Range("B6").Value = "US Dollar"
$wsIn.Range("B2").Copy()
$wsIn.Range("B6").PasteSpecial(-4122)

# 4. Clear the stray empty D column cells (D1, D2)
$wsIn.Range("D1").Clear()
$wsIn.Range("D2").Clear()

# 5. Clear the stray empty C11 cell
$wsIn.Range("C11").Clear()

# Mirror the same product-name fix on the output sheet
$wsOut.Range("B1").Value = "775-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-Late Repayment"

# Update selections/active cells to match final workbook view state
$wsIn.Range("B2:B3").Select()
$wsOut.Range("B1").Select()
